# The deck currently carries two embedded themes:
#   ppt/theme/theme1.xml -> "Integral"      (active theme, used by the slide master)
#   ppt/theme/theme2.xml -> "Office Theme"  (only used by the notes master)
#
# The author switched the presentation's applied design from "Integral" to
# the default "Office Theme" colour palette. Reproduce that by pushing the
# Office Theme's 12 theme colours into the presentation's active colour
# scheme (this is exactly what Design > Variants > Colors does in the UI).
#
# Colour order for ThemeColorScheme.Item(n) is the standard OOXML clrScheme
# slot order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
# RGB values are passed the same way VBA's RGB()/OLE_COLOR does: r + g*256 + b*65536.

$p = $ppt.ActivePresentation
$slides = $p.Slides.Range()
$colors = $slides.ThemeColorScheme

$colors.Item(1).RGB  = 0        # dk1      000000
$colors.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388  # dk2      44546A
$colors.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501  # accent2  ED7D31
$colors.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Item(8).RGB  = 49407    # accent4  FFC000
$colors.Item(9).RGB  = 12874308 # accent5  4472C4
$colors.Item(10).RGB = 4697456  # accent6  70AD47
$colors.Item(11).RGB = 12673797 # hlink    0563C1
$colors.Item(12).RGB = 7491477  # folHlink 954F72
